$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Attendance list: mark "C" (presente) for the 02/05/2022 class (columns M = manha,
# N = tarde) for every student row that was already marked present on 29/04 (column L),
# the same way attendance was recorded for every earlier class date. Row 25 has no
# attendance recorded yet for that student, so it is left untouched.
for ($r = 3; $r -le 49; $r++) {
    if ($r -eq 25) { continue }
    $rangeAddr = "M" + $r + ":N" + $r
    $ws.Range($rangeAddr).Value = "C"
}

# Leave the selection where the author ended up after filling in the new column.
$ws.Range("N51").Select()
